$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The backing Power Query ("Query1") was refreshed, producing fewer rows and
# updated "Days remaining" figures. Mirror that refreshed result set:
#  - the "BNT323-01" trial row is gone
#  - the "REDEFINE HF" trial row is gone
#  - "REJOICE (MK-5909-003)" days remaining dropped from 22 to 16
#  - "REMASTER (CLOU)" days remaining dropped from 42 to 36

# Row 3 is "BNT323-01" - remove it, shifting everything below up.
$ws.Rows.Item(3).Delete()

# After that shift, the former "REDEFINE HF" row (originally row 10) is now
# row 9 - remove it too.
$ws.Rows.Item(9).Delete()

# Update the two changed "Days remaining" values.
$ws.Cells.Item(8, 2).Value = 16
$ws.Cells.Item(10, 2).Value = 36

# Keep the "ExternalData_1" defined name (the Power Query's cached range)
# in sync with the now-smaller result set.
$wb.Names.Item("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$C`$10"
